$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.030.09"
$ws.Range("E2").Value = "  +0.51%  "
$ws.Range("D3").Value = "1.659.41"
$ws.Range("E3").Value = "  +2.14%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.59%  "
$ws.Range("D5").Value = "'310.10"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "  -0.43%  "
$ws.Range("D7").Value = "'0.3904"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.58%  "
$ws.Range("D9").Value = "'51.31"
$ws.Range("E9").Value = "  +2.81%  "
$ws.Range("D10").Value = "'1.370"
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").Value = "'0.08507"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("D13").Value = "'24.02"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "'7.216"
$ws.Range("E14").Value = "  +3.37%  "
$ws.Range("D15").Value = "'8.036"
$ws.Range("E15").Value = "  +8.35%  "
$ws.Range("E16").Value = "  +3.50%  "
$ws.Range("D17").Value = "1.656.68"
$ws.Range("E17").Value = "  +1.71%  "
$ws.Range("D18").Value = "'94.65"
$ws.Range("E18").Value = "  +1.92%  "
$ws.Range("D19").Value = "'0.06992"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "'20.01"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "'6.993"
$ws.Range("E21").Value = "  +2.21%  "
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("D23").Value = "'13.70"
$ws.Range("E23").Value = "  +2.40%  "
$ws.Range("D24").Value = "24.029.94"
$ws.Range("E24").Value = "  +0.46%  "
$ws.Range("D25").Value = "'2.493"
$ws.Range("E25").Value = "  +4.46%  "
$ws.Range("D26").Value = "'3.109"
$ws.Range("E26").Value = "  +8.57%  "
$ws.Range("D27").Value = "'22.27"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("D28").Value = "'153.93"
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("D29").Value = "'140.36"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").Value = "'5.313"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "'7.891"
$ws.Range("E31").Value = "  +1.88%  "
$ws.Range("E32").Value = "  +2.34%  "
$ws.Range("D33").Value = "1.835.38"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "'1.056"
$ws.Range("E34").Value = "  +11.24%  "
$ws.Range("D35").Value = "'0.08174"
$ws.Range("E35").Value = "  +2.84%  "
$ws.Range("D36").Value = "'0.03002"
$ws.Range("E36").Value = "  +4.17%  "
$ws.Range("D37").Value = "'11.15"
$ws.Range("E37").Value = "  +9.13%  "
$ws.Range("D38").Value = "'6.722"
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("D39").Value = "'0.2712"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").Value = "'0.09150"
$ws.Range("E40").Value = "  -0.03%  "
$ws.Range("D41").Value = "'13.65"
$ws.Range("E41").Value = "  +4.17%  "
$ws.Range("D42").Value = "'0.7594"
$ws.Range("E42").Value = "  +2.17%  "
$ws.Range("D43").Value = "'1.425"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("D44").Value = "'16.65"
$ws.Range("E44").Value = "  +4.76%  "
$ws.Range("D45").Value = "'0.7041"
$ws.Range("E45").Value = "  +3.37%  "
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("D47").Value = "'4.098"
$ws.Range("E47").Value = "  +0.40%  "
$ws.Range("E48").Value = "  -0.42%  "
$ws.Range("D49").Value = "'0.08303"
$ws.Range("E49").Value = "  +0.75%  "
$ws.Range("D50").Value = "'135.66"
$ws.Range("E50").Value = "  +1.93%  "
$ws.Range("D51").Value = "'1.241"
$ws.Range("E51").Value = "  -0.52%  "
